$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Column widths (A & B get noticeably wider to fit new content)
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 74
$ws.Columns.Item(2).ColumnWidth = 48

# ---------------------------------------------------------------
# 2. Row 4 - the old "login" row: drop its hyperlink/shrink its height
#    (the hyperlink is being moved down to the new "regenerate-otp" row)
# ---------------------------------------------------------------
$ws.Cells.Item(4, 1).Hyperlinks.Delete()
$ws.Rows.Item(4).RowHeight = 201.6

# ---------------------------------------------------------------
# 3. Row 5 - verify-account endpoint
# ---------------------------------------------------------------
$ws.Range("A5").Value = "http://localhost:8080/kinMel/verify-account?email=suman.yhhits@gmail.com&otp=927583"
$ws.Range("A5").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

$ws.Range("B5").Value = "PUT"
$ws.Range("B5").VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignBottom

$ws.Range("D5").Value = "OTP verified.Now you can login"
$ws.Range("D5").WrapText = $true

# ---------------------------------------------------------------
# 4. Row 6 - regenerate-otp endpoint (receives the hyperlink)
# ---------------------------------------------------------------
$ws.Range("A6").Value = "http://localhost:8080/kinMel/regenerate-otp?email=dssuman222@gmail.com"
$ws.Hyperlinks.Add($ws.Range("A6"), "http://localhost:8080/auth/login")
$ws.Range("A6").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

$ws.Range("B6").Value = "PUT"
$ws.Range("B6").VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignBottom

# ---------------------------------------------------------------
# 5. Rows 7 & 8 - blank centered placeholder cells under column A
# ---------------------------------------------------------------
$ws.Range("A7").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$ws.Range("A8").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

# ---------------------------------------------------------------
# 6. Row 9 - blank centered placeholder cell under column C
# ---------------------------------------------------------------
$ws.Range("C9").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

# ---------------------------------------------------------------
# 7. Selection / view state
# ---------------------------------------------------------------
$ws.Range("D6").Select()
